# Scheduled-runner style refresh of cached market-board figures across the
# item-level "Profits" sheets (one tab per crafting job). Each tab has the
# same layout: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ. These are plain cached numbers (no formulas in the
# workbook), so the refresh is just a series of direct value writes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 968179.6
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 1113271.5
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 3339814.5
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -3340150.5
$ws.Range("H19").Value = 1640.7307
$ws.Range("I19").Value = 1335.8572
$ws.Range("J19").Value = 1996.4166
$ws.Range("K19").Value = 1335.8572
$ws.Range("L19").Value = 1996.4166
$ws.Range("M19").Value = -1160.8572
$ws.Range("N19").Value = -2346.4166
$ws.Range("H40").Value = 1113626.8
$ws.Range("I40").Value = 1668773.1
$ws.Range("J40").Value = 3334
$ws.Range("K40").Value = 1668773.1
$ws.Range("L40").Value = 3334
$ws.Range("M40").Value = -1668598.1
$ws.Range("N40").Value = -3684
$ws.Range("H44").Value = 4447.5
$ws.Range("I44").Value = 4447.5
$ws.Range("K44").Value = 4447.5
$ws.Range("M44").Value = -3985.5
$ws.Range("H113").Value = 7057.4165
$ws.Range("I113").Value = 4049
$ws.Range("J113").Value = 7659.1
$ws.Range("K113").Value = 4049
$ws.Range("L113").Value = 7659.1
$ws.Range("M113").Value = -795
$ws.Range("N113").Value = -14167.1
$ws.Range("H132").Value = 3164.1052
$ws.Range("I132").Value = 2715.3333
$ws.Range("J132").Value = 3933.4285
$ws.Range("K132").Value = 8145.999899999999
$ws.Range("L132").Value = 11800.2855
$ws.Range("M132").Value = -5615.999899999999
$ws.Range("N132").Value = -16860.2855
$ws.Range("H137").Value = 3209.2632
$ws.Range("J137").Value = 2160.2856
$ws.Range("L137").Value = 6480.8568
$ws.Range("N137").Value = -11580.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 815.46155
$ws.Range("I2").Value = 660.2
$ws.Range("J2").Value = 1333
$ws.Range("K2").Value = 660.2
$ws.Range("L2").Value = 1333
$ws.Range("M2").Value = -547.2
$ws.Range("N2").Value = -1559
$ws.Range("H32").Value = 8011.061
$ws.Range("J32").Value = 26093.9
$ws.Range("L32").Value = 26093.9
$ws.Range("N32").Value = -26667.9
$ws.Range("H45").Value = 2513
$ws.Range("I45").Value = 2012
$ws.Range("J45").Value = 3014
$ws.Range("K45").Value = 2012
$ws.Range("L45").Value = 3014
$ws.Range("M45").Value = -1635
$ws.Range("N45").Value = -3768
$ws.Range("H116").Value = 815.46155
$ws.Range("I116").Value = 660.2
$ws.Range("J116").Value = 1333
$ws.Range("K116").Value = 660.2
$ws.Range("L116").Value = 1333
$ws.Range("M116").Value = 1633.8
$ws.Range("N116").Value = -5921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 734
$ws.Range("J3").Value = 1333
$ws.Range("L3").Value = 1333
$ws.Range("N3").Value = -1561
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = ""
$ws.Range("N28").Value = 0
$ws.Range("H52").Value = 35797.5
$ws.Range("J52").Value = 35797.5
$ws.Range("L52").Value = 35797.5
$ws.Range("N52").Value = -36323.5
$ws.Range("H117").Value = 59999
$ws.Range("J117").Value = 59999
$ws.Range("L117").Value = 59999
$ws.Range("N117").Value = -69177
$ws.Range("H121").Value = 35797.5
$ws.Range("J121").Value = 35797.5
$ws.Range("L121").Value = 35797.5
$ws.Range("N121").Value = -39291.5
$ws.Range("H141").Value = 49996
$ws.Range("J141").Value = 49996
$ws.Range("L141").Value = 49996
$ws.Range("N141").Value = -60356

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7045.933
$ws.Range("I22").Value = 7520.7144
$ws.Range("J22").Value = 399
$ws.Range("K22").Value = 7520.7144
$ws.Range("L22").Value = 399
$ws.Range("M22").Value = -7170.7144
$ws.Range("N22").Value = -1099
$ws.Range("H52").Value = 98891.42999999999
$ws.Range("I52").Value = 97998
$ws.Range("J52").Value = 99040.336
$ws.Range("K52").Value = 97998
$ws.Range("L52").Value = 99040.336
$ws.Range("M52").Value = -97704
$ws.Range("N52").Value = -99628.336
$ws.Range("H99").Value = 5575.375
$ws.Range("I99").Value = 6164.75
$ws.Range("J99").Value = 2628.5
$ws.Range("K99").Value = 6164.75
$ws.Range("L99").Value = 2628.5
$ws.Range("M99").Value = -4666.75
$ws.Range("N99").Value = -5624.5
$ws.Range("H126").Value = 5575.375
$ws.Range("I126").Value = 6164.75
$ws.Range("J126").Value = 2628.5
$ws.Range("K126").Value = 18494.25
$ws.Range("L126").Value = 7885.5
$ws.Range("M126").Value = -16024.25
$ws.Range("N126").Value = -12825.5
$ws.Range("H131").Value = 69757.836
$ws.Range("J131").Value = 69757.836
$ws.Range("L131").Value = 69757.836
$ws.Range("N131").Value = -79837.836
$ws.Range("H134").Value = 1181.1111
$ws.Range("I134").Value = 1181.1111
$ws.Range("K134").Value = 3543.3333
$ws.Range("M134").Value = -1008.3333
$ws.Range("H139").Value = 125000
$ws.Range("J139").Value = 125000
$ws.Range("L139").Value = 125000
$ws.Range("N139").Value = -135280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 198997.67
$ws.Range("J37").Value = 198997.67
$ws.Range("L37").Value = 596993.01
$ws.Range("N37").Value = -597217.01
$ws.Range("H87").Value = 2001
$ws.Range("I87").Value = 2001
$ws.Range("K87").Value = 6003
$ws.Range("M87").Value = -4755
$ws.Range("H90").Value = 2001
$ws.Range("I90").Value = 2001
$ws.Range("K90").Value = 18009
$ws.Range("M90").Value = -11769
$ws.Range("H92").Value = 397
$ws.Range("I92").Value = 398
$ws.Range("J92").Value = 396.66666
$ws.Range("K92").Value = 1194
$ws.Range("L92").Value = 1189.99998
$ws.Range("M92").Value = 54
$ws.Range("N92").Value = -3685.99998
$ws.Range("H113").Value = 2797
$ws.Range("J113").Value = 4007.375
$ws.Range("L113").Value = 12022.125
$ws.Range("N113").Value = -16362.125
$ws.Range("H121").Value = 1464.091
$ws.Range("I121").Value = 178
$ws.Range("K121").Value = 534
$ws.Range("M121").Value = 776
$ws.Range("H131").Value = 38064.29
$ws.Range("J131").Value = 5070.852
$ws.Range("L131").Value = 15212.556
$ws.Range("N131").Value = -25292.556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4998.4287
$ws.Range("I80").Value = 4998
$ws.Range("J80").Value = 4999.5
$ws.Range("K80").Value = 4998
$ws.Range("L80").Value = 4999.5
$ws.Range("M80").Value = -4000
$ws.Range("N80").Value = -6995.5
$ws.Range("H83").Value = 4998.4287
$ws.Range("I83").Value = 4998
$ws.Range("J83").Value = 4999.5
$ws.Range("K83").Value = 24990
$ws.Range("L83").Value = 24997.5
$ws.Range("M83").Value = -19998
$ws.Range("N83").Value = -34981.5
$ws.Range("H132").Value = 4500.778
$ws.Range("I132").Value = 6674.6665
$ws.Range("K132").Value = 20023.9995
$ws.Range("M132").Value = -17493.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4095.1667
$ws.Range("I7").Value = 4149
$ws.Range("J7").Value = 4056.7144
$ws.Range("K7").Value = 4149
$ws.Range("L7").Value = 4056.7144
$ws.Range("M7").Value = -4037
$ws.Range("N7").Value = -4280.7144
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15452
$ws.Range("H63").Value = 69333
$ws.Range("I63").Value = 69500
$ws.Range("J63").Value = 68999
$ws.Range("K63").Value = 69500
$ws.Range("L63").Value = 68999
$ws.Range("M63").Value = -68751
$ws.Range("N63").Value = -70497
$ws.Range("H66").Value = 69333
$ws.Range("I66").Value = 69500
$ws.Range("J66").Value = 68999
$ws.Range("K66").Value = 208500
$ws.Range("L66").Value = 206997
$ws.Range("M66").Value = -204756
$ws.Range("N66").Value = -214485
$ws.Range("H100").Value = 3916.8
$ws.Range("I100").Value = 2594
$ws.Range("K100").Value = 2594
$ws.Range("M100").Value = -2053
$ws.Range("H126").Value = 4095.1667
$ws.Range("I126").Value = 4149
$ws.Range("J126").Value = 4056.7144
$ws.Range("K126").Value = 12447
$ws.Range("L126").Value = 12170.1432
$ws.Range("M126").Value = -9977
$ws.Range("N126").Value = -17110.1432
$ws.Range("H136").Value = 4399.6
$ws.Range("I136").Value = 1499.5
$ws.Range("K136").Value = 4498.5
$ws.Range("M136").Value = -1948.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1762.129
$ws.Range("I136").Value = 1638.1333
$ws.Range("K136").Value = 4914.3999
$ws.Range("M136").Value = -2364.3999
